$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Carry row 20 (last row of the existing block) forward to a bottom-border
# style, matching the other closed dialogue blocks in the sheet, by copying
# the formatting from a row that already has that "closing" style (row 16).
$ws.Range("A16:E16").Copy() | Out-Null
$ws.Range("A20:E20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New block: rows 21-23 ---------------------------------------------
# (literal backslash-n / backslash-apostrophe sequences, matching the
# source script's own escaping convention used throughout this sheet)
# Row 21
$ws.Range("C21").Value = " Team [CS:X]Raider[CR] is a legendary\nexploration team."
$ws.Range("A21").Value = "SCRIPT/T01P02A/us0409.ssb"
# Row 22
$ws.Range("C22").Value = " The team consists of three\nmembers: [CS:N]Gallade[CR], [CS:N]Roserade[CR], and [CS:N]Rhyperior[CR]."
# Row 23
$ws.Range("C23").Value = " They\'re said to never be denied\nfrom getting the treasures they target."

# Russian translations
$ws.Range("D21").Value = " Команда [CS:X]Рейдер[CR] это легендарная\nкоманда исследователей."
$ws.Range("D22").Value = " В ней состоят трое: [CS:N]Галлейд[CR],\n[CS:N]Роузрэйд[CR] и [CS:N]Райпериор[CR]."
$ws.Range("D23").Value = " Говорят, что если они начнут\nохоту за каким-либо сокровищем, их ничто\nне остановит."

# "Converted"/garbled strings
$ws.Range("E21").Value = " Ëïíàîäà [CS:X]Ñåêäåñ[CR] üóï ìåãåîäàñîàÿ\nëïíàîäà éòòìåäïâàóåìåê."
$ws.Range("E22").Value = " Â îåê òïòóïÿó óñïå: [CS:N]Ãàììåêä[CR],\n[CS:N]Ñïôèñüêä[CR] é [CS:N]Ñàêðåñéïñ[CR]."
$ws.Range("E23").Value = " Ãïâïñÿó, œóï åòìé ïîé îàœîôó\nïöïóô èà ëàëéí-ìéáï òïëñïâéþåí, éö îéœóï\nîå ïòóàîïâéó."

# Row 22 filename (entered last, after the other row-22/23 content)
$ws.Range("A22").Value = "SCRIPT/T01P02A/us2005.ssb"

# Line numbers
$ws.Range("B21").Value = 41
$ws.Range("B22").Value = 44
$ws.Range("B23").Value = 47

# --- Formatting for the new rows: same look as the other open (non-bordered)
# rows, so copy from row 19 (A/B/C/D/E all style 4/5) which is the most
# recent row written in the matching style.
$ws.Range("A19:E19").Copy() | Out-Null
$ws.Range("A21:E21").PasteSpecial(-4122) | Out-Null
$ws.Range("A22:E22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B13:E13").Copy() | Out-Null
$ws.Range("B23:E23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row heights ---------------------------------------------------------
$ws.Rows.Item(21).RowHeight = 43.2
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Rows.Item(23).RowHeight = 31.8

# --- Selection / view -----------------------------------------------------
$ws.Range("D19").Select()
